$d = $word.ActiveDocument

# --- Change 1: Insert a new paragraph ("There are 3 separate jupyter notebook
#     files...") right after the "For the 3 tasks of this project..." intro
#     paragraph, i.e. right before the two blank bold paragraphs that lead
#     into the "Company Name extraction" heading.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "For the 3 tasks of this project*") {
        $anchor = $p
        break
    }
}

if ($anchor -ne $null) {
    $anchor.Range.InsertParagraphAfter()
    $newPara = $anchor.Next()
    $newPara.Range.Text = "There are 3 separate jupyter notebook files for the 3 different tasks attached in the github repo."
}

# --- Change 2: Tweak wording in the "Sample Preparation" paragraph:
#     "... on the right side as the snippet I looked at to generate features."
#     becomes "... on the right as the snippet to generate features."
$d.Content.Find.Execute(
    "on the left and 50 characters on the right side as the snippet I looked at to generate features. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "on the left and 50 characters on the right as the snippet to generate features. ", 2)
